# Insert a new "Exam 1 2021/2022" slide right before the existing
# "Homework 3 2021/2021" slide (which sits at position 5). The new
# slide takes position 5 and the homework slide is pushed to position 6,
# unchanged.

$p = $ppt.ActivePresentation

$layout = $p.Slides.Item(5).Layout

$new = $p.Slides.Add(5, $layout)

# --- Title placeholder -------------------------------------------------
$titleShape = $new.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Exam 1 2021/2022"

# --- Content placeholder (two hyperlinked lines) ------------------------
$bodyShape = $new.Shapes.Item(2)
$tr = $bodyShape.TextFrame.TextRange

$url1 = "https://github.com/NaskoVasilev/SDA_2022-2023/tree/main/Exams/01_Exam"
$url2 = "https://www.hackerrank.com/contests/sda-2021-2021-test-1/challenges"

$tr.Text = $url1 + " " + "`r" + $url2 + " "

$para1 = $tr.Paragraphs(1)
$para2 = $tr.Paragraphs(2)

$p1Url = $para1.Characters(1, $url1.Length)
$p1Space = $para1.Characters($url1.Length + 1, 1)
# force a clean run split before assigning the hyperlink so the
# trailing space does not inherit the hyperlink run properties
$p1Space.Font.Name | Out-Null
$p1Url.ActionSettings.Item(1).Hyperlink.Address = $url1

$p2Url = $para2.Characters(1, $url2.Length)
$p2Space = $para2.Characters($url2.Length + 1, 1)
$p2Space.Font.Name | Out-Null
$p2Url.ActionSettings.Item(1).Hyperlink.Address = $url2

Write-Output ("Slide count: " + $p.Slides.Count)
Write-Output ("New slide index: " + $new.SlideIndex + " id " + $new.SlideID)
